$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 84: copy the date formatting from A83 (reuses existing style, avoids
# creating a new cellXf), then set the new values.
$ws.Range("A83").Copy()
$ws.Range("A84").PasteSpecial(-4122)
$ws.Range("A84").Value = 43818
$ws.Range("B84").Value = 2203.4753172042001
$ws.Range("C84").Value = 2207.0300000000002
$ws.Range("D84").Formula = "=100*(B84-C84)/C84"
$ws.Range("E84").Value = 169
$ws.Range("F84").Value = "New CRM opened 12/11/2020"

# Row 85: only column F populated.
$ws.Range("F85").Value = "New CRM opened 12/11/2021"

# Update the active selection to reflect the new last cell, same as the
# author's workbook did after editing.
$null = $ws.Range("A85").Select()
